$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(15, 17, 18)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "The American Journal of Gastroenterology"
    $ws.Range("G$r").Value = "https://openalex.org/S66441642"
    $ws.Range("H$r").Value = "Lippincott Williams & Wilkins"
    $ws.Range("I$r").Value = "0002-9270"
    $ws.Range("V$r").Value = "'FALSE"
}
